$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking values are preserved exactly
# (matches the original inline-string cell type rather than being reinterpreted as numbers).
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '42.575.88'
$ws.Range("E2").Value = '  -1.70%  '

$ws.Range("D3").Value = '2.516.10'
$ws.Range("E3").Value = '  -3.32%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").Value = '307.75'
$ws.Range("E5").Value = '  -2.91%  '

$ws.Range("D6").Value = '100.57'
$ws.Range("E6").Value = '  +2.50%  '

$ws.Range("D7").Value = '0.569'
$ws.Range("E7").Value = '  -1.78%  '

$ws.Range("E8").Value = '  +0.23%  '

$ws.Range("D9").Value = '0.524'
$ws.Range("E9").Value = '  -3.38%  '

$ws.Range("D10").Value = '35.91'
$ws.Range("E10").Value = '  -0.57%  '

$ws.Range("D11").Value = '0.0800'
$ws.Range("E11").Value = '  -2.10%  '

$ws.Range("D12").Value = '7.28'
$ws.Range("E12").Value = '  -3.97%  '

$ws.Range("E13").Value = '  -0.10%  '

$ws.Range("D14").Value = '2.918.25'
$ws.Range("E14").Value = '  -2.86%  '

$ws.Range("D15").Value = '15.56'
$ws.Range("E15").Value = '  +1.90%  '

$ws.Range("D16").Value = '2.521.90'
$ws.Range("E16").Value = '  -2.85%  '

$ws.Range("D17").Value = '0.804'
$ws.Range("E17").Value = '  -5.52%  '

$ws.Range("D18").Value = '42.578.43'

$ws.Range("D19").Value = '6.68'
$ws.Range("E19").Value = '  -2.86%  '

$ws.Range("D20").Value = '0.0₃0945'
$ws.Range("E20").Value = '  -2.76%  '

$ws.Range("D21").Value = '12.09'
$ws.Range("E21").Value = '  -5.39%  '

$ws.Range("D22").Value = '69.32'
$ws.Range("E22").Value = '  -0.46%  '

$ws.Range("D23").Value = '243.27'
$ws.Range("E23").Value = '  -4.62%  '

$ws.Range("D24").Value = '2.87'
$ws.Range("E24").Value = '  -3.69%  '

$ws.Range("D25").Value = '2.03'
$ws.Range("E25").Value = '  -2.86%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("D27").Value = '25.82'
$ws.Range("E27").Value = '  -5.43%  '

$ws.Range("E28").Value = '  -3.90%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '10.10'
$ws.Range("E29").Value = '  -2.39%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '38.73'
$ws.Range("E30").Value = '  -6.84%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '155.74'
$ws.Range("E31").Value = '  -0.62%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '5.74'
$ws.Range("E32").Value = '  -2.59%  '

$ws.Range("D33").Value = '2.76'
$ws.Range("E33").Value = '  +10.01%  '

$ws.Range("D34").Value = '0.0785'
$ws.Range("E34").Value = '  -3.47%  '

$ws.Range("D35").Value = '2.62'
$ws.Range("E35").Value = '  -2.68%  '

$ws.Range("D36").Value = '2.01'
$ws.Range("E36").Value = '  -7.47%  '

$ws.Range("D37").Value = '3.17'
$ws.Range("E37").Value = '  -9.43%  '

$ws.Range("D38").Value = '18.09'
$ws.Range("E38").Value = '  -3.99%  '

$ws.Range("D39").Value = '0.111'
$ws.Range("E39").Value = '  -1.70%  '

$ws.Range("D40").Value = '0.118'
$ws.Range("E40").Value = '  -0.43%  '

$ws.Range("D41").Value = '4.25'
$ws.Range("E41").Value = '  +5.51%  '

$ws.Range("D42").Value = '22.09'
$ws.Range("E42").Value = '  -3.51%  '

$ws.Range("E43").Value = '  +0.16%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0298'
$ws.Range("E44").Value = '  -2.43%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '3.28'
$ws.Range("E45").Value = '  +0.35%  '

$ws.Range("D46").Value = '1.978.99'
$ws.Range("E46").Value = '  -1.86%  '

$ws.Range("D47").Value = '8.81'
$ws.Range("E47").Value = '  -2.33%  '

$ws.Range("D48").Value = '2.772.76'
$ws.Range("E48").Value = '  -3.01%  '

$ws.Range("D49").Value = '79.89'
$ws.Range("E49").Value = '  -4.73%  '

$ws.Range("D50").Value = '0.189'
$ws.Range("E50").Value = '  -3.57%  '

$ws.Range("D51").Value = '72.14'
$ws.Range("E51").Value = '  -4.06%  '

# Restore column D to the default (unstyled) appearance now that values are set as text.
$dRange.Style = "Normal"
